# ============================================================
# edit.ps1
# Applies the "stubble sim" rerun:
#  1. Renames existing sheet tabs (positions 4-8) to their new
#     names, in an order that avoids transient name collisions.
#  2. Appends four new sheets ("f", "i", "k", "v") at the end.
#  3. Rewrites the A1:D5 data block on every sheet (1-12) with
#     the freshly-simulated values from the rerun.
# ============================================================

function Set-SheetData {
    param($ws, $data)
    for ($r = 0; $r -lt $data.Length; $r++) {
        $row = $data[$r]
        for ($c = 0; $c -lt $row.Length; $c++) {
            $ws.Cells.Item($r + 1, $c + 1).Value = [double]$row[$c]
        }
    }
}

$data1 = @(
    @("0.8489756835642335","0","0","0"),
    @("0.1166744211042993","0.7393627540995942","0.5167543779040673","0.009258364106652164"),
    @("0.009937343463686236","0.07533564712301603","0.1389725936762536","0.2387538065204591"),
    @("0.01005900884216587","0.07625075290223962","0.1405837976771043","0.2381096899657622"),
    @("0.01435354302561505","0.1090508458751502","0.2036892307425748","0.5138781394071265")
)

$data2 = @(
    @("0.8210040926182817","0.1136333442085682","0","0"),
    @("0.1072076018592264","0.4870638938672912","0.3035086869258043","0.004153653652502913"),
    @("0.03610662388276849","0.2008247596554859","0.3502079867871424","0.256145448951541"),
    @("0.035521196631256","0.1975707186406024","0.3445534101313601","0.2521845648611721"),
    @("0.000160485008467316","0.000907283628052448","0.001729916155693289","0.4875163325347838")
)

$data3 = @(
    @("0.7632431705833963","0.01461359346040558","0","0"),
    @("0.1776586271990664","0.7281681340483785","0.6955126932651354","0.184988845145639"),
    @("0","0","0","0"),
    @("0","0","0","0"),
    @("0.05909820221753723","0.2572182724912158","0.3044873067348645","0.8150111548543609")
)

$data4 = @(
    @("0","0","0","0"),
    @("0.9903741887482146","0.4841387322429547","0","0"),
    @("0.0047950222870563","0.2569722396608778","0.4981421473610299","0.4981421473610304"),
    @("0.004830787359102568","0.2588889420483039","0.5018576858346971","0.5018576858346965"),
    @("1.60562644740132e-09","8.604786371092232e-08","1.668042729493081e-07","1.668042729493079e-07")
)

$data5 = @(
    @("0.9420548328005345","0.6029168001815662","0.2740467478864517","0"),
    @("0.008078472041071923","0.02747651264889068","0.1142795874827693","0.1391482107355601"),
    @("0","0","0","0"),
    @("2.764120770839165e-18","9.034397199979058e-18","3.460311420382062e-17","1.363577806854834e-17"),
    @("0.0498666951583937","0.1696066871695426","0.7054236646307785","0.858932403467895")
)

$data6 = @(
    @("0.4271062341794111","0","0","0"),
    @("0.5695200870426755","0.9744378431396262","0","0"),
    @("7.869531232707829e-13","5.955550572052952e-12","1.771096990464771e-10","1.028671089127642e-11"),
    @("0.00250788657463974","0.01899745592075483","0.6983594029114579","0.1099537198267961"),
    @("0.0008657922024865962","0.006564700933663254","0.3016405969114322","0.8900462801629173")
)

$data7 = @(
    @("0.4271062341794111","0","0","0"),
    @("0.5695200870426755","0.9744378431396262","0","0"),
    @("7.869531232707829e-13","5.955550572052952e-12","1.771096990464771e-10","1.028671089127642e-11"),
    @("0.00250788657463974","0.01899745592075483","0.6983594029114579","0.1099537198267961"),
    @("0.0008657922024865962","0.006564700933663254","0.3016405969114322","0.8900462801629173")
)

$data8 = @(
    @("0.797811283493076","0.04431981106641727","0","0"),
    @("0.2008885663752003","0.9456236137619409","0.513481222667956","0"),
    @("1.384126259710869e-14","1.779888970714261e-15","0","0"),
    @("0.0007123345067536939","0.00550928280829495","0.2655597553020279","0.3476765817938068"),
    @("0.000587815624956144","0.004547292363345072","0.2209590220300156","0.6523234182061932")
)

$dataTemplate = @(
    @("0.797811283493076","0.04431981106641727","0","0"),
    @("0.2008885663752003","0.9456236137619409","0.513481222667956","0"),
    @("1.384126259710869e-14","1.779888970714261e-15","0","0"),
    @("0.0007123345067536939","0.00550928280829495","0.2655597553020279","0.3476765817938068"),
    @("0.000587815624956144","0.004547292363345072","0.2209590220300156","0.6523234182061932")
)
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------
# Step 1: rename the existing tabs (positions 4-8) to their
# new names. Renaming is done tail-first along the rename
# chain  f -> z -> h -> l -> r -> of  so that no intermediate
# step ever requires two sheets to share the same name.
# ------------------------------------------------------------
$wb.Worksheets.Item(5).Name = "of"   # was "r"
$wb.Worksheets.Item(7).Name = "r"    # was "l"
$wb.Worksheets.Item(8).Name = "l"    # was "h"
$wb.Worksheets.Item(4).Name = "h"    # was "z"
$wb.Worksheets.Item(6).Name = "z"    # was "f"

# ------------------------------------------------------------
# Step 2: add four brand-new input sheets at the end of the
# workbook: "f", "i", "k", "v".
# ------------------------------------------------------------
$newNames = @("f", "i", "k", "v")
foreach ($name in $newNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $name
}

# ------------------------------------------------------------
# Step 3: write the new simulated A1:D5 values into every sheet
# (by position, which is stable across the rename operations).
# ------------------------------------------------------------
Set-SheetData $wb.Worksheets.Item(1)  $data1         # w
Set-SheetData $wb.Worksheets.Item(2)  $data2         # b
Set-SheetData $wb.Worksheets.Item(3)  $data3         # o
Set-SheetData $wb.Worksheets.Item(4)  $data4         # h
Set-SheetData $wb.Worksheets.Item(5)  $data5         # of
Set-SheetData $wb.Worksheets.Item(6)  $data6         # z
Set-SheetData $wb.Worksheets.Item(7)  $data7         # r
Set-SheetData $wb.Worksheets.Item(8)  $data8         # l
Set-SheetData $wb.Worksheets.Item(9)  $dataTemplate  # f
Set-SheetData $wb.Worksheets.Item(10) $dataTemplate  # i
Set-SheetData $wb.Worksheets.Item(11) $dataTemplate  # k
Set-SheetData $wb.Worksheets.Item(12) $dataTemplate  # v

Write-Host "Sheets now:"
foreach ($ws in $wb.Worksheets) {
    Write-Host " -" $ws.Name
}
